$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before C (current energy column) to hold "artist",
# shifting energy..popularity from C:J to D:K.
$ws.Columns("C").Insert()

# Header row (left to right): song, artist, energy, dance, liveness, valence,
# tempo, instrumental, acoustic, popularity, genres
$ws.Range("C1").Value = "artist"
$ws.Range("L1").Value = "genres"

# Match the header formatting (bold, centered, bordered) used by the other
# header cells for the new "genres" header cell.
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Artist column (C) values for each song row, top to bottom
$ws.Range("C2").Value = "The Weeknd"
$ws.Range("C3").Value = "Drake"
$ws.Range("C4").Value = "Drake"
$ws.Range("C5").Value = "Coldplay"
$ws.Range("C7").Value = "Drake"
$ws.Range("C8").Value = "Mike Posner"
$ws.Range("C9").Value = "Justin Bieber"
$ws.Range("C10").Value = "Ariana Grande"
$ws.Range("C11").Value = "Nicki Minaj"

# Genres column (L) values for each song row, top to bottom
$ws.Range("L2").Value = "['canadian pop', 'pop']"
$ws.Range("L3").Value = "['canadian pop', 'hip hop', 'pop rap', 'rap']"
$ws.Range("L4").Value = "['canadian pop', 'hip hop', 'pop rap', 'rap']"
$ws.Range("L5").Value = "['permanent wave', 'pop', 'pop christmas', 'rock']"
$ws.Range("L7").Value = "['canadian pop', 'hip hop', 'pop rap', 'rap']"
$ws.Range("L8").Value = "['dance pop', 'pop', 'pop rap', 'post-teen pop', 'tropical house']"
$ws.Range("L9").Value = "['canadian pop', 'dance pop', 'pop', 'pop christmas', 'post-teen pop']"
$ws.Range("L10").Value = "['dance pop', 'pop', 'pop christmas', 'post-teen pop']"
$ws.Range("L11").Value = "['dance pop', 'dwn trap', 'hip pop', 'pop', 'pop rap']"
